$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> [new A value (date-as-number), new E value]
$data = @{
    3 = @(20150926, 5)
    4 = @(20150927, 9)
    5 = @(20150928, 14)
    6 = @(20150929, 8)
    7 = @(20150930, 15)
    8 = @(20150931, 5)
    9 = @(20150932, 12)
    10 = @(20150933, 5)
    11 = @(20150934, 18)
    12 = @(20150935, 5)
    13 = @(20150936, 14)
    14 = @(20150937, 17)
    15 = @(20150938, 19)
    16 = @(20150939, 20)
    17 = @(20150940, 20)
    18 = @(20150941, 13)
    19 = @(20150942, 8)
    20 = @(20150943, 6)
    21 = @(20150944, 8)
    22 = @(20150945, 6)
    23 = @(20150946, 8)
    24 = @(20150947, 6)
    25 = @(20150948, 12)
    26 = @(20150949, 11)
    27 = @(20150950, 11)
    28 = @(20150951, 15)
    29 = @(20150952, 18)
    30 = @(20150953, 20)
    31 = @(20150954, 20)
    32 = @(20150955, 8)
    33 = @(20150956, 20)
    34 = @(20150957, 5)
    35 = @(20150958, 11)
    36 = @(20150959, 5)
    37 = @(20150960, 17)
    38 = @(20150961, 15)
    39 = @(20150962, 8)
    40 = @(20150963, 10)
    41 = @(20150964, 14)
    42 = @(20150965, 7)
    43 = @(20150966, 14)
    44 = @(20150967, 8)
    45 = @(20150968, 19)
    46 = @(20150969, 16)
    47 = @(20150970, 16)
    48 = @(20150971, 13)
    49 = @(20150972, 19)
    50 = @(20150973, 5)
    51 = @(20150974, 14)
    52 = @(20150975, 9)
    53 = @(20150976, 13)
    54 = @(20150977, 11)
    55 = @(20150978, 20)
    56 = @(20150979, 20)
    57 = @(20150980, 11)
    58 = @(20150981, 5)
    59 = @(20150982, 15)
    60 = @(20150983, 10)
    61 = @(20150984, 7)
    62 = @(20150985, 9)
    63 = @(20150986, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
